# Apply trade #24 close update to the live trading results workbook.
$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# Summary sheet: refresh aggregate metrics after the new closed trade.
# ----------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1199.62   # Current Capital
$wsSummary.Range("B4").Value = -0.38     # Total P&L $
$wsSummary.Range("B5").Value = -0.32     # Total P&L %
$wsSummary.Range("B6").Value = 24        # Total Trades
$wsSummary.Range("B7").Value = 7         # Winning Trades
$wsSummary.Range("B9").Value = 29.17     # Win Rate %

# ----------------------------------------------------------------------
# Strategy Status sheet: update the MarketMaking strategy row (row 4).
# ----------------------------------------------------------------------
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C4").Value = 99.62      # Capital
$wsStatus.Range("D4").Value = 24         # Trades
$wsStatus.Range("E4").Value = -0.38      # P&L $
$wsStatus.Range("F4").Value = -0.38      # P&L %
$wsStatus.Range("G4").Value = 29.17      # Win Rate %

# ----------------------------------------------------------------------
# All Trades & MarketMaking sheets: append the newly closed trade #24
# as row 25.
# ----------------------------------------------------------------------
$sheetNames = @("All Trades", "MarketMaking")
foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Cells.Item(25, 1).Value = 24
    # Leading apostrophe forces text so the date-like string isn't
    # auto-converted into a date serial number (matches the source data,
    # which stores dates/times as literal text).
    $ws.Cells.Item(25, 2).Value = "'2026-02-17"
    $ws.Cells.Item(25, 3).Value = "08:02:59"
    $ws.Cells.Item(25, 4).Value = "MarketMaking"
    $ws.Cells.Item(25, 5).Value = "DOWN"
    $ws.Cells.Item(25, 6).Value = 0.58
    $ws.Cells.Item(25, 7).Value = 0.59
    $ws.Cells.Item(25, 8).Value = "CLOSED"
    $ws.Cells.Item(25, 9).Value = 1.7241
    $ws.Cells.Item(25, 10).Value = 0.01
    $ws.Cells.Item(25, 11).Value = 99.62
    $ws.Cells.Item(25, 12).Value = 0
    $ws.Cells.Item(25, 13).Value = 0
    $ws.Cells.Item(25, 14).Value = 0.6
    $ws.Cells.Item(25, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item(25, 16).Value = "early_exit"
    $ws.Cells.Item(25, 17).Value = 0.14
}
